$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (shifts the old "Session 6" row and everything
# below it down by one row), making room for the new "Session 6" entry
# (Matt Bombyk / IPUMS talk) that now precedes the renamed "Session 7".
$ws.Rows("4:4").Insert()
$ws.Rows("4:4").RowHeight = 45

# New row 4: Session 6 - Linking Administrative Data: The IPUMS Experience
$ws.Range("A4").Value = "1:15PM"
$ws.Range("B4").Value = "1:35PM"
$ws.Range("C4").Value = "20 minutes"
$ws.Range("D4").Value = "Session 6"
$ws.Range("E4").Value = "Linking Administrative Data: The IPUMS Experience"
$ws.Range("F4").Value = "[Matt Bombyk](https://dataifa.github.io/difa-project/comingsoon.html)"

# D4/F4 pick up the same (time) number format as the rest of the row.
$ws.Range("D4").NumberFormat = $ws.Range("A4").NumberFormat
$ws.Range("F4").NumberFormat = $ws.Range("A4").NumberFormat

# Row 5 (previously row 4): "Session 6 (Research Presentations)" becomes
# "Session 7 (Research Presentations)" and its times shift later.
$ws.Range("A5").Value = "1:35PM"
$ws.Range("B5").Value = "2:35PM"
$ws.Range("D5").Value = "Session 7 (Research Presentations)"

# Row 6 (previously row 5): Break - times shift later.
$ws.Range("A6").Value = "2:35PM"
$ws.Range("B6").Value = "2:50PM"

# Row 7 (previously row 6): Activity 1 - times shift later.
$ws.Range("A7").Value = "2:50PM"
$ws.Range("B7").Value = "3:30PM"

# Row 8 (previously row 7): Activity 2 - times shift later.
$ws.Range("A8").Value = "3:30PM"
$ws.Range("B8").Value = "4:10PM"

# Row 9 (previously row 8): Activity 3 - times shift later.
$ws.Range("A9").Value = "4:10PM"
$ws.Range("B9").Value = "4:50PM"

# Row 10 (previously row 9): Day 2 Wrap-up - start time shifts later.
$ws.Range("A10").Value = "4:50PM"

# Update selection to match the final saved state.
$ws.Range("D10").Select()
